# Apply Hill Climbing experiment data edits (mirrors the authored OOXML diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 56: I56 becomes a literal number (was a shared-string "207.01") ---
$ws.Cells.Item(56, 9).Value = 207.01

# Row 57
$ws.Cells.Item(57, 3).Value = 210.96
$ws.Cells.Item(57, 4).Value = 193
$ws.Cells.Item(57, 5).Value = 27
$ws.Cells.Item(57, 6).Value = 141.6
$ws.Cells.Item(57, 7).Value = 447
$ws.Cells.Item(57, 8).Value = 83
$ws.Cells.Item(57, 9).Value = 148.79
$ws.Cells.Item(57, 10).Value = 242
$ws.Cells.Item(57, 11).Value = 23.82
$ws.Cells.Item(57, 12).Value = 15.33
$ws.Cells.Item(57, 13).Value = 15.53

# Row 58
$ws.Cells.Item(58, 2).Value = 224
$ws.Cells.Item(58, 3).Value = 214.19
$ws.Cells.Item(58, 4).Value = 200
$ws.Cells.Item(58, 5).Value = 76
$ws.Cells.Item(58, 6).Value = 339.18
$ws.Cells.Item(58, 7).Value = 881
$ws.Cells.Item(58, 8).Value = 77
$ws.Cells.Item(58, 9).Value = 127.53
$ws.Cells.Item(58, 10).Value = 203
$ws.Cells.Item(58, 11).Value = 22.98
$ws.Cells.Item(58, 12).Value = 16.85
$ws.Cells.Item(58, 13).Value = 15.57

# Row 59
$ws.Cells.Item(59, 2).Value = 231
$ws.Cells.Item(59, 3).Value = 218.33
$ws.Cells.Item(59, 4).Value = 206
$ws.Cells.Item(59, 5).Value = 149
$ws.Cells.Item(59, 6).Value = 479.82
$ws.Cells.Item(59, 7).Value = 1019
$ws.Cells.Item(59, 8).Value = 40
$ws.Cells.Item(59, 9).Value = 105.03
$ws.Cells.Item(59, 10).Value = 166
$ws.Cells.Item(59, 11).Value = 23.7
$ws.Cells.Item(59, 12).Value = 16.97
$ws.Cells.Item(59, 13).Value = 15.87

# Row 60
$ws.Cells.Item(60, 2).Value = 231
$ws.Cells.Item(60, 3).Value = 218.53
$ws.Cells.Item(60, 4).Value = 209
$ws.Cells.Item(60, 5).Value = 137
$ws.Cells.Item(60, 6).Value = 595.29
$ws.Cells.Item(60, 7).Value = 1471
$ws.Cells.Item(60, 8).Value = 42
$ws.Cells.Item(60, 9).Value = 103.68
$ws.Cells.Item(60, 10).Value = 157
$ws.Cells.Item(60, 11).Value = 23.55
$ws.Cells.Item(60, 12).Value = 17.07
$ws.Cells.Item(60, 13).Value = 15.92

# Row 61
$ws.Cells.Item(61, 2).Value = 234
$ws.Cells.Item(61, 3).Value = 219.98
$ws.Cells.Item(61, 4).Value = 209
$ws.Cells.Item(61, 5).Value = 443
$ws.Cells.Item(61, 6).Value = 1310.59
$ws.Cells.Item(61, 7).Value = 2441
$ws.Cells.Item(61, 8).Value = 18
$ws.Cells.Item(61, 9).Value = 95.71
$ws.Cells.Item(61, 10).Value = 150
$ws.Cells.Item(61, 11).Value = 23.57
$ws.Cells.Item(61, 12).Value = 17.13
$ws.Cells.Item(61, 13).Value = 16.15

# Row 63
$ws.Cells.Item(63, 1).Value = 'Only A, B and C packages used; "infinite" supply of packages; Objective function: Maximize Tot. Value'

# Row 64
$ws.Cells.Item(64, 1).Value = 'Random Mutation: 10 package; Rotations: enabled'

# Row 65
$ws.Cells.Item(65, 1).Value = 'Neighbourhood size'
$ws.Cells.Item(65, 2).Value = 'Best value'
$ws.Cells.Item(65, 3).Value = 'Average value'
$ws.Cells.Item(65, 4).Value = 'Worst value'
$ws.Cells.Item(65, 5).Value = 'Best runtime (ms)'
$ws.Cells.Item(65, 6).Value = 'Average runtime (ms)'
$ws.Cells.Item(65, 7).Value = 'Worst runtime (ms)'
$ws.Cells.Item(65, 8).Value = 'Least gaps'
$ws.Cells.Item(65, 9).Value = 'Average gaps'
$ws.Cells.Item(65, 10).Value = 'Most gaps'
$ws.Cells.Item(65, 11).Value = 'Average A'
$ws.Cells.Item(65, 12).Value = 'Average B'
$ws.Cells.Item(65, 13).Value = 'Average C'

# Row 66
$ws.Cells.Item(66, 1).Value = 10
$ws.Cells.Item(66, 2).Value = 217
$ws.Cells.Item(66, 3).Value = 200.93
$ws.Cells.Item(66, 4).Value = 184
$ws.Cells.Item(66, 5).Value = 2
$ws.Cells.Item(66, 6).Value = 21.11
$ws.Cells.Item(66, 7).Value = 191
$ws.Cells.Item(66, 8).Value = 122
$ws.Cells.Item(66, 9).Value = 202.22
$ws.Cells.Item(66, 10).Value = 289
$ws.Cells.Item(66, 11).Value = 22.45
$ws.Cells.Item(66, 12).Value = 15.52
$ws.Cells.Item(66, 13).Value = 14.3

# Row 67
$ws.Cells.Item(67, 1).Value = 50
$ws.Cells.Item(67, 2).Value = 227
$ws.Cells.Item(67, 3).Value = 212.45
$ws.Cells.Item(67, 4).Value = 194
$ws.Cells.Item(67, 5).Value = 15
$ws.Cells.Item(67, 6).Value = 137.72
$ws.Cells.Item(67, 7).Value = 525
$ws.Cells.Item(67, 8).Value = 71
$ws.Cells.Item(67, 9).Value = 141.3
$ws.Cells.Item(67, 10).Value = 234
$ws.Cells.Item(67, 11).Value = 24.21
$ws.Cells.Item(67, 12).Value = 15.13
$ws.Cells.Item(67, 13).Value = 15.86

# Row 68
$ws.Cells.Item(68, 1).Value = 100
$ws.Cells.Item(68, 2).Value = 229
$ws.Cells.Item(68, 3).Value = 216.22
$ws.Cells.Item(68, 4).Value = 197
$ws.Cells.Item(68, 5).Value = 117
$ws.Cells.Item(68, 6).Value = 349.26
$ws.Cells.Item(68, 7).Value = 1022
$ws.Cells.Item(68, 8).Value = 58
$ws.Cells.Item(68, 9).Value = 119.17
$ws.Cells.Item(68, 10).Value = 202
$ws.Cells.Item(68, 11).Value = 24.95
$ws.Cells.Item(68, 12).Value = 15.93
$ws.Cells.Item(68, 13).Value = 15.53

# Row 69
$ws.Cells.Item(69, 1).Value = 200
$ws.Cells.Item(69, 2).Value = 231
$ws.Cells.Item(69, 3).Value = 219.92
$ws.Cells.Item(69, 4).Value = 208
$ws.Cells.Item(69, 5).Value = 310
$ws.Cells.Item(69, 6).Value = 701.31
$ws.Cells.Item(69, 7).Value = 1485
$ws.Cells.Item(69, 8).Value = 30
$ws.Cells.Item(69, 9).Value = '''98.6'
$ws.Cells.Item(69, 10).Value = 157
$ws.Cells.Item(69, 11).Value = 25
$ws.Cells.Item(69, 12).Value = 16.18
$ws.Cells.Item(69, 13).Value = 16.04

# Row 70
$ws.Cells.Item(70, 1).Value = 250
$ws.Cells.Item(70, 2).Value = 235
$ws.Cells.Item(70, 3).Value = 220.4
$ws.Cells.Item(70, 4).Value = 205
$ws.Cells.Item(70, 5).Value = 256
$ws.Cells.Item(70, 6).Value = 929.8
$ws.Cells.Item(70, 7).Value = 1791
$ws.Cells.Item(70, 8).Value = 29
$ws.Cells.Item(70, 9).Value = 95.17
$ws.Cells.Item(70, 10).Value = 172
$ws.Cells.Item(70, 11).Value = 24.77
$ws.Cells.Item(70, 12).Value = 16.51
$ws.Cells.Item(70, 13).Value = 16.01

# Row 71
$ws.Cells.Item(71, 1).Value = 500
$ws.Cells.Item(71, 2).Value = 235
$ws.Cells.Item(71, 3).Value = 222.45
$ws.Cells.Item(71, 4).Value = 212
$ws.Cells.Item(71, 5).Value = 478
$ws.Cells.Item(71, 6).Value = 1676.42
$ws.Cells.Item(71, 7).Value = 2925
$ws.Cells.Item(71, 8).Value = 34
$ws.Cells.Item(71, 9).Value = 84.7
$ws.Cells.Item(71, 10).Value = 150
$ws.Cells.Item(71, 11).Value = 24.77
$ws.Cells.Item(71, 12).Value = 16.26
$ws.Cells.Item(71, 13).Value = 16.62

# Row 73
$ws.Cells.Item(73, 1).Value = 'Only A, B and C packages used; "infinite" supply of packages; Objective function: Maximize Tot. Value'

# Row 74
$ws.Cells.Item(74, 1).Value = 'Random Mutation: 20 package; Rotations: enabled'

# Row 75
$ws.Cells.Item(75, 1).Value = 'Neighbourhood size'
$ws.Cells.Item(75, 2).Value = 'Best value'
$ws.Cells.Item(75, 3).Value = 'Average value'
$ws.Cells.Item(75, 4).Value = 'Worst value'
$ws.Cells.Item(75, 5).Value = 'Best runtime (ms)'
$ws.Cells.Item(75, 6).Value = 'Average runtime (ms)'
$ws.Cells.Item(75, 7).Value = 'Worst runtime (ms)'
$ws.Cells.Item(75, 8).Value = 'Least gaps'
$ws.Cells.Item(75, 9).Value = 'Average gaps'
$ws.Cells.Item(75, 10).Value = 'Most gaps'
$ws.Cells.Item(75, 11).Value = 'Average A'
$ws.Cells.Item(75, 12).Value = 'Average B'
$ws.Cells.Item(75, 13).Value = 'Average C'

# Row 76
$ws.Cells.Item(76, 1).Value = 10
$ws.Cells.Item(76, 2).Value = 218
$ws.Cells.Item(76, 3).Value = 199.86
$ws.Cells.Item(76, 4).Value = 185
$ws.Cells.Item(76, 5).Value = 7
$ws.Cells.Item(76, 6).Value = 32.84
$ws.Cells.Item(76, 7).Value = 159
$ws.Cells.Item(76, 8).Value = 110
$ws.Cells.Item(76, 9).Value = 206.82
$ws.Cells.Item(76, 10).Value = 279
$ws.Cells.Item(76, 11).Value = 20.64
$ws.Cells.Item(76, 12).Value = 15.86
$ws.Cells.Item(76, 13).Value = 14.9

# Row 77
$ws.Cells.Item(77, 1).Value = 50
$ws.Cells.Item(77, 2).Value = 220
$ws.Cells.Item(77, 3).Value = 209.18
$ws.Cells.Item(77, 4).Value = 194
$ws.Cells.Item(77, 5).Value = 32
$ws.Cells.Item(77, 6).Value = 215.77
$ws.Cells.Item(77, 7).Value = 637
$ws.Cells.Item(77, 8).Value = 99
$ws.Cells.Item(77, 9).Value = 158.94
$ws.Cells.Item(77, 10).Value = 235
$ws.Cells.Item(77, 11).Value = 23.52
$ws.Cells.Item(77, 12).Value = 15.08
$ws.Cells.Item(77, 13).Value = 15.66

# Row 78
$ws.Cells.Item(78, 1).Value = 100
$ws.Cells.Item(78, 2).Value = 226
$ws.Cells.Item(78, 3).Value = 212.74
$ws.Cells.Item(78, 4).Value = 197
$ws.Cells.Item(78, 5).Value = 60
$ws.Cells.Item(78, 6).Value = 313.79
$ws.Cells.Item(78, 7).Value = 645
$ws.Cells.Item(78, 8).Value = 78
$ws.Cells.Item(78, 9).Value = 140.68
$ws.Cells.Item(78, 10).Value = 217
$ws.Cells.Item(78, 11).Value = 24.38
$ws.Cells.Item(78, 12).Value = 14.75
$ws.Cells.Item(78, 13).Value = 16.12

# Row 79
$ws.Cells.Item(79, 1).Value = 200
$ws.Cells.Item(79, 2).Value = 228
$ws.Cells.Item(79, 3).Value = 215.72
$ws.Cells.Item(79, 4).Value = 202
$ws.Cells.Item(79, 5).Value = 214
$ws.Cells.Item(79, 6).Value = 675.24
$ws.Cells.Item(79, 7).Value = 1226
$ws.Cells.Item(79, 8).Value = 68
$ws.Cells.Item(79, 9).Value = 125.85
$ws.Cells.Item(79, 10).Value = 200
$ws.Cells.Item(79, 11).Value = 25.29
$ws.Cells.Item(79, 12).Value = 14.3
$ws.Cells.Item(79, 13).Value = 16.53

# Row 80
$ws.Cells.Item(80, 1).Value = 250
$ws.Cells.Item(80, 2).Value = 228
$ws.Cells.Item(80, 3).Value = 216.36
$ws.Cells.Item(80, 4).Value = 206
$ws.Cells.Item(80, 5).Value = 285
$ws.Cells.Item(80, 6).Value = 917.7
$ws.Cells.Item(80, 7).Value = 2095
$ws.Cells.Item(80, 8).Value = 66
$ws.Cells.Item(80, 9).Value = 121.42
$ws.Cells.Item(80, 10).Value = 175
$ws.Cells.Item(80, 11).Value = 25.34
$ws.Cells.Item(80, 12).Value = 14.71
$ws.Cells.Item(80, 13).Value = '''16.3'

# Row 81
$ws.Cells.Item(81, 1).Value = 500
$ws.Cells.Item(81, 2).Value = 232
$ws.Cells.Item(81, 3).Value = 218.54
$ws.Cells.Item(81, 4).Value = 203
$ws.Cells.Item(81, 5).Value = 803
$ws.Cells.Item(81, 6).Value = 2739.65
$ws.Cells.Item(81, 7).Value = 4735
$ws.Cells.Item(81, 8).Value = 33
$ws.Cells.Item(81, 9).Value = 108.9
$ws.Cells.Item(81, 10).Value = 185
$ws.Cells.Item(81, 11).Value = 25.32
$ws.Cells.Item(81, 12).Value = 15.02
$ws.Cells.Item(81, 13).Value = 16.5

# --- View: selection + scroll position (topLeftCell is not persisted by this host) ---
try {
    $excel.ActiveWindow.ScrollRow = 51
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("K82").Select()
